{"js": "// Localize the FAQ copy from Indian-Rupee / \"EMI\" phrasing to American\n// banking phrasing (\"$\" instead of \"Rs.\", \"monthly payment\" instead of\n// \"EMI\"), per the commit message \"changes data to reflect american banks\n// and names for HCL workshop\".\n//\n// Each entry below is a unique, unambiguous substring of the document's\n// text (verified against the original) together with its replacement.\nconst replacements = [\n  // \"Why I am not able to select loan amount below 25000?\"\n  [\"loan amount below 25000?\", \"loan amount below $25000?\"],\n\n  // \"Minimum loan amount for this product is 25000.\"\n  [\"product is 25000.\", \"product is $25000.\"],\n\n  // \"Under this scheme maximum loan is Rs.1,00,000/- & maxi...\"\n  [\"maximum loan is Rs.1,00,000/-\", \"maximum loan is $1,00,000/-\"],\n\n  // \"Can I select my EMI amount?\" (heading)\n  [\"select my EMI amount?\", \"select my monthly payment amount?\"],\n\n  // \"...loan amount & tenure to suit your EMI \" (body, before Lumpsum para)\n  [\"suit your EMI \", \"suit your monthly payment \"],\n\n  // \"This letter is sent to the Borrowers with loan of more than Rs. 75,000/-\"\n  [\"more than Rs. 75,000/-\", \"more than $75,000/-\"],\n\n  // \"...if the eligibility is 1,00,000/- can I take...\"\n  [\"eligibility is 1,00,000/-\", \"eligibility is $1,00,000/-\"],\n\n  // \"...2 loans of Rs.50,000/-, one now & other Next month)?\"\n  [\"2 loans of Rs.50,000/-\", \"2 loans of $50,000/-\"],\n];\n\nfor (const [searchText, replacementText] of replacements) {\n  const results = context.document.body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + searchText);\n  }\n\n  // Replace every match (each search string is unique in this document,\n  // so this normally fires exactly once).\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacementText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Localize the FAQ copy from Indian-Rupee / \"EMI\" phrasing to American\n# banking phrasing (\"$\" instead of \"Rs.\", \"monthly payment\" instead of\n# \"EMI\"), per the commit message \"changes data to reflect american banks\n# and names for HCL workshop\".\n#\n# Each pair below is a unique, unambiguous substring of the document's\n# text (verified against the original) together with its replacement.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"loan amount below 25000?\"; Replace = \"loan amount below `$25000?\" },\n    @{ Find = \"product is 25000.\";        Replace = \"product is `$25000.\" },\n    @{ Find = \"maximum loan is Rs.1,00,000/-\"; Replace = \"maximum loan is `$1,00,000/-\" },\n    @{ Find = \"select my EMI amount?\";    Replace = \"select my monthly payment amount?\" },\n    @{ Find = \"suit your EMI \";           Replace = \"suit your monthly payment \" },\n    @{ Find = \"more than Rs. 75,000/-\";   Replace = \"more than `$75,000/-\" },\n    @{ Find = \"eligibility is 1,00,000/-\"; Replace = \"eligibility is `$1,00,000/-\" },\n    @{ Find = \"2 loans of Rs.50,000/-\";   Replace = \"2 loans of `$50,000/-\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $found = $range.Find.Execute($pair.Find, $true, $false, $false, $false, $false, $true, 1, $false, $pair.Replace, 2)\n    if (-not $found) {\n        throw \"No match found for: $($pair.Find)\"\n    }\n}\n"}
